# Generate Report for Handoff
# Updates the generated handoff/handback artifact names (old GUID -> new GUID,
# old content-hash -> new content-hash) and the associated timestamps across
# the "Overview", "zh-cn" and "de-de" sheets, keeping each sheet's existing
# hyperlink target URL intact but refreshing its displayed text.

$wb = $excel.ActiveWorkbook

$oldBase = "e34a4148-e28c-4206-ab7c-41882eb065c1"
$newBase = "e1a71351-3ca9-4e52-a4cf-3d9aa822e7a4"

$oldHash = "0e02e62ee839cedf304266e7f8c86730166f8548"
$newHash = "cebe26dfc03180e89d7c64921921cd2b181e439b"

# ---------------------------------------------------------------------------
# Sheet "Overview"  (File Name / Path And Name / ... / Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newBase.md"
$wsOverview.Range("G2").Value = "2016-08-24 09:00:23"

# B2 carries a hyperlink; keep its (unchanged) target address but refresh the
# displayed text to reference the new file name.
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbb9039b047224b6f685194929ae14ab36ab70a2/e2e/$oldBase.md",
    $null,
    $null,
    "e2e\$newBase.md"
)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("G2").Value = "$newBase.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-24 08:59:58"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbb9039b047224b6f685194929ae14ab36ab70a2/e2e/$oldBase.md",
    $null,
    $null,
    "$newBase.md"
)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("G2").Value = "$newBase.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-24 09:00:23"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A2"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dbb9039b047224b6f685194929ae14ab36ab70a2/e2e/$oldBase.md",
    $null,
    $null,
    "$newBase.md"
)
